$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new "Path" column (B) with the source file paths for each
# media item now that the ListView click updates the MediaPlayer's source.
# Values are entered in the same order they were authored (video items
# first, then image items) so new shared-string entries land in the same
# order as the original edit.
$ws.Range("B4").Value = "C:/Users/Sagar/Downloads/vid.mp4"
$ws.Range("B6").Value = "C:/Users/Sagar/Downloads/production ID_4713259.mp4"
$ws.Range("B7").Value = "C:/Users/Sagar/Downloads/istockphoto.mp4"
$ws.Range("B2").Value = "C:/Users/Sagar/Downloads/swan.jpg"
$ws.Range("B3").Value = "C:/Users/Sagar/Downloads/scene.jpeg"
$ws.Range("B5").Value = "C:/Users/Sagar/Downloads/photo.jpeg"
$ws.Range("B8").Value = "C:/Users/Sagar/Downloads/images.jpeg"

# Move the active selection to B9, matching the last interaction location.
$ws.Range("B9").Select()
